$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain display strings scraped from the site.
# Most look like text (thousand-separator dots, e.g. "28.729.83") but a few
# are bare decimals (e.g. "154.40") that Excel would otherwise auto-convert
# to a number and silently drop the trailing zero / formatting. Force those
# specific cells to Text format before writing so the literal string sticks.
$ws.Range("D2").Value = "28.729.83"
$ws.Range("E2").Value = "  -1.50%  "

$ws.Range("D3").Value = "1.806.08"
$ws.Range("E3").Value = "  -0.98%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.12"
$ws.Range("E5").Value = "  -1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5953"
$ws.Range("E6").Value = "  -2.33%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2788"
$ws.Range("E8").Value = "  -0.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06851"
$ws.Range("E9").Value = "  -3.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.38"
$ws.Range("E10").Value = "  -0.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07495"
$ws.Range("E11").Value = "  -2.15%  "

$ws.Range("D12").Value = "1.808.08"
$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6253"
$ws.Range("E14").Value = "  -1.04%  "

$ws.Range("D15").Value = "2.051.78"
$ws.Range("E15").Value = "  -0.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009347"
$ws.Range("E16").Value = "  -6.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.92"
$ws.Range("E17").Value = "  -3.39%  "

$ws.Range("D18").Value = "28.662.62"
$ws.Range("E18").Value = "  -1.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.514"
$ws.Range("E19").Value = "  -5.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "211.65"
$ws.Range("E21").Value = "  -6.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.49"
$ws.Range("E22").Value = "  -2.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.878"
$ws.Range("E23").Value = "  -1.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.40"
$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.896"
$ws.Range("E26").Value = "  -1.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1277"
$ws.Range("E27").Value = "  -2.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.45"
$ws.Range("E28").Value = "  -0.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.425"
$ws.Range("E29").Value = "  -4.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06243"
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.424"
$ws.Range("E31").Value = "  -1.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.795"
$ws.Range("E32").Value = "  -0.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.767"
$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.728"
$ws.Range("E34").Value = "  -0.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.069"
$ws.Range("E35").Value = "  -4.82%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.493"
$ws.Range("E37").Value = "  -2.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.714"
$ws.Range("E38").Value = "  -0.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.505"
$ws.Range("E39").Value = "  -0.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01718"
$ws.Range("E40").Value = "  -1.10%  "

$ws.Range("D41").Value = "1.140.45"
$ws.Range("E41").Value = "  -6.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8788"
$ws.Range("E42").Value = "  -3.08%  "

$ws.Range("E43").Value = "  +0.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.53"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").Value = "1.964.33"
$ws.Range("E45").Value = "  -0.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.63"
$ws.Range("E46").Value = "  -3.24%  "

$ws.Range("E47").Value = "  -3.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.609"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.415"
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("E50").Value = "  -0.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4485"
$ws.Range("E51").Value = "  -1.65%  "
